$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header row from Spanish to English
$ws.Range("A1").Value = "Module"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Requirement ID"
$ws.Range("D1").Value = "Requirement"
$ws.Range("E1").Value = "Considerations"

# Update the active selection to E2 (as reflected in the saved view state)
$ws.Range("E2").Select()
